$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated simulation results values (columns B, D, E, F, G, K, L, N) for rows 2-25
$ws.Range("B2").Value = 1.246748183887064 ; $ws.Range("D2").Value = 0.100861413049131 ; $ws.Range("E2").Value = 0.1123289386138282 ; $ws.Range("F2").Value = 2.146450807788369 ; $ws.Range("G2").Value = 0.002528076292829612 ; $ws.Range("K2").Value = 1.051048226445658 ; $ws.Range("L2").Value = 0.2473636141730253 ; $ws.Range("N2").Value = 2.950801149231197
$ws.Range("B3").Value = 1.219779472695791 ; $ws.Range("D3").Value = 0.1013388829962043 ; $ws.Range("E3").Value = 0.1118091654777409 ; $ws.Range("F3").Value = 2.109131284187768 ; $ws.Range("G3").Value = 0.002533022867710667 ; $ws.Range("K3").Value = 0.9632996346812206 ; $ws.Range("L3").Value = 0.2344631084123847 ; $ws.Range("N3").Value = 2.955732365344502
$ws.Range("B4").Value = 1.203904745414206 ; $ws.Range("D4").Value = 0.1016488144719325 ; $ws.Range("E4").Value = 0.1114849743318507 ; $ws.Range("F4").Value = 2.087328611008772 ; $ws.Range("G4").Value = 0.002536220283943275 ; $ws.Range("K4").Value = 0.91002312539797 ; $ws.Range("L4").Value = 0.226699275821062 ; $ws.Range("N4").Value = 2.959471811802956
$ws.Range("B5").Value = 1.197607745780971 ; $ws.Range("D5").Value = 0.1017793225283228 ; $ws.Range("E5").Value = 0.1113515652256996 ; $ws.Range("F5").Value = 2.078722417702167 ; $ws.Range("G5").Value = 0.002537563676452535 ; $ws.Range("K5").Value = 0.8884624480444643 ; $ws.Range("L5").Value = 0.2235747983106933 ; $ws.Range("N5").Value = 2.961174249550922
$ws.Range("B6").Value = 1.196572528415487 ; $ws.Range("D6").Value = 0.1018012472793366 ; $ws.Range("E6").Value = 0.1113293334705294 ; $ws.Range("F6").Value = 2.077310163066286 ; $ws.Range("G6").Value = 0.00253778919080865 ; $ws.Range("K6").Value = 0.8848913315631535 ; $ws.Range("L6").Value = 0.2230583529439372 ; $ws.Range("N6").Value = 2.961467713272654
$ws.Range("B7").Value = 1.203819125071448 ; $ws.Range("D7").Value = 0.1016505575198856 ; $ws.Range("E7").Value = 0.1114831804245306 ; $ws.Range("F7").Value = 2.087211418231135 ; $ws.Range("G7").Value = 0.002536238237453509 ; $ws.Range("K7").Value = 0.9097317452149696 ; $ws.Range("L7").Value = 0.2266569789233728 ; $ws.Range("N7").Value = 2.959494048910898
$ws.Range("B8").Value = 1.237307383169849 ; $ws.Range("D8").Value = 0.1010225614622939 ; $ws.Range("E8").Value = 0.112150750866026 ; $ws.Range("F8").Value = 2.133351710239822 ; $ws.Range("G8").Value = 0.002529748704113816 ; $ws.Range("K8").Value = 1.020666850504341 ; $ws.Range("L8").Value = 0.2428828032791586 ; $ws.Range("N8").Value = 2.952353501945495
$ws.Range("B9").Value = 1.308409982705911 ; $ws.Range("D9").Value = 0.0999243646223924 ; $ws.Range("E9").Value = 0.1134211668923655 ; $ws.Range("F9").Value = 2.232704125096276 ; $ws.Range("G9").Value = 0.002518287669319223 ; $ws.Range("K9").Value = 1.24305448750232 ; $ws.Range("L9").Value = 0.2759577867196015 ; $ws.Range("N9").Value = 2.944015416127485
$ws.Range("B10").Value = 1.363973635046278 ; $ws.Range("D10").Value = 0.09919917919837928 ; $ws.Range("E10").Value = 0.1143329307094501 ; $ws.Range("F10").Value = 2.31118772450111 ; $ws.Range("G10").Value = 0.002510629698982776 ; $ws.Range("K10").Value = 1.409513179812336 ; $ws.Range("L10").Value = 0.3010397651094934 ; $ws.Range("N10").Value = 2.941368945131671
$ws.Range("B11").Value = 1.389976272287186 ; $ws.Range("D11").Value = 0.0988870951990144 ; $ws.Range("E11").Value = 0.1147434606068503 ; $ws.Range("F11").Value = 2.348102262297886 ; $ws.Range("G11").Value = 0.002507309597348315 ; $ws.Range("K11").Value = 1.485932920811877 ; $ws.Range("L11").Value = 0.3126236262877029 ; $ws.Range("N11").Value = 2.940926409058378
$ws.Range("B12").Value = 1.399927386029134 ; $ws.Range("D12").Value = 0.09877148643744249 ; $ws.Range("E12").Value = 0.1148983434340769 ; $ws.Range("F12").Value = 2.362256415614041 ; $ws.Range("G12").Value = 0.002506075737758935 ; $ws.Range("K12").Value = 1.514973105534182 ; $ws.Range("L12").Value = 0.3170353864800575 ; $ws.Range("N12").Value = 2.940868780068683
$ws.Range("B13").Value = 1.397779588548019 ; $ws.Range("D13").Value = 0.09879627035562066 ; $ws.Range("E13").Value = 0.1148650117262715 ; $ws.Range("F13").Value = 2.359200245278004 ; $ws.Range("G13").Value = 0.00250634043322648 ; $ws.Range("K13").Value = 1.508714235740342 ; $ws.Range("L13").Value = 0.3160841109900616 ; $ws.Range("N13").Value = 2.940876293933343
$ws.Range("B14").Value = 1.390792862362559 ; $ws.Range("D14").Value = 0.09887753241663688 ; $ws.Range("E14").Value = 0.1147562142604119 ; $ws.Range("F14").Value = 2.349263210754259 ; $ws.Range("G14").Value = 0.002507207618897629 ; $ws.Range("K14").Value = 1.488320027748841 ; $ws.Range("L14").Value = 0.3129860782373584 ; $ws.Range("N14").Value = 2.940919461628795
$ws.Range("B15").Value = 1.38652689774338 ; $ws.Range("D15").Value = 0.09892764285918076 ; $ws.Range("E15").Value = 0.1146894986634823 ; $ws.Range("F15").Value = 2.343199369077411 ; $ws.Range("G15").Value = 0.002507741837519127 ; $ws.Range("K15").Value = 1.475841273556853 ; $ws.Range("L15").Value = 0.3110917316573847 ; $ws.Range("N15").Value = 2.940960235784345
$ws.Range("B16").Value = 1.362288920861516 ; $ws.Range("D16").Value = 0.09921993379196437 ; $ws.Range("E16").Value = 0.1143060191998186 ; $ws.Range("F16").Value = 2.308799757451595 ; $ws.Range("G16").Value = 0.00251084995534809 ; $ws.Range("K16").Value = 1.404533118466418 ; $ws.Range("L16").Value = 0.300286248934114 ; $ws.Range("N16").Value = 2.941413226178199
$ws.Range("B17").Value = 1.347605778360872 ; $ws.Range("D17").Value = 0.09940381258101993 ; $ws.Range("E17").Value = 0.1140697054536135 ; $ws.Range("F17").Value = 2.28800792737627 ; $ws.Range("G17").Value = 0.002512798483268153 ; $ws.Range("K17").Value = 1.360967400907214 ; $ws.Range("L17").Value = 0.293702120849332 ; $ws.Range("N17").Value = 2.941886453550907
$ws.Range("B18").Value = 1.339228798924978 ; $ws.Range("D18").Value = 0.09951125020266716 ; $ws.Range("E18").Value = 0.1139333832987228 ; $ws.Range("F18").Value = 2.276163022315046 ; $ws.Range("G18").Value = 0.002513934625992146 ; $ws.Range("K18").Value = 1.335975052941308 ; $ws.Range("L18").Value = 0.2899314653012368 ; $ws.Range("N18").Value = 2.942230275605297
$ws.Range("B19").Value = 1.336404240851436 ; $ws.Range("D19").Value = 0.09954791418450171 ; $ws.Range("E19").Value = 0.1138871571309412 ; $ws.Range("F19").Value = 2.272172087052894 ; $ws.Range("G19").Value = 0.002514321953716188 ; $ws.Range("K19").Value = 1.32752427836698 ; $ws.Range("L19").Value = 0.2886575903513915 ; $ws.Range("N19").Value = 2.942358976925476
$ws.Range("B20").Value = 1.349161747229459 ; $ws.Range("D20").Value = 0.0993840649003559 ; $ws.Range("E20").Value = 0.1140949026891001 ; $ws.Range("F20").Value = 2.290209445851445 ; $ws.Range("G20").Value = 0.002512589466001795 ; $ws.Range("K20").Value = 1.365598256195881 ; $ws.Range("L20").Value = 0.2944013178738913 ; $ws.Range("N20").Value = 2.941828660503916
$ws.Range("B21").Value = 1.392842197496094 ; $ws.Range("D21").Value = 0.09885359393021886 ; $ws.Range("E21").Value = 0.1147881860766002 ; $ws.Range("F21").Value = 2.352177188668009 ; $ws.Range("G21").Value = 0.002506952271457392 ; $ws.Range("K21").Value = 1.494307529307662 ; $ws.Range("L21").Value = 0.3138953600069101 ; $ws.Range("N21").Value = 2.940903794463566
$ws.Range("B22").Value = 1.421998895429624 ; $ws.Range("D22").Value = 0.09852188928304528 ; $ws.Range("E22").Value = 0.115237945699751 ; $ws.Range("F22").Value = 2.393699708925624 ; $ws.Range("G22").Value = 0.002503404321274093 ; $ws.Range("K22").Value = 1.579019896030445 ; $ws.Range("L22").Value = 0.3267828388030978 ; $ws.Range("N22").Value = 2.940940434327672
$ws.Range("B23").Value = 1.406381686230674 ; $ws.Range("D23").Value = 0.09869755138964109 ; $ws.Range("E23").Value = 0.1149981951115242 ; $ws.Range("F23").Value = 2.371444370598965 ; $ws.Range("G23").Value = 0.002505285501078989 ; $ws.Range("K23").Value = 1.533752518212395 ; $ws.Range("L23").Value = 0.3198910349104267 ; $ws.Range("N23").Value = 2.940862062664038
$ws.Range("B24").Value = 1.348458092625947 ; $ws.Range("D24").Value = 0.09939298745874936 ; $ws.Range("E24").Value = 0.1140835124625692 ; $ws.Range("F24").Value = 2.289213801166511 ; $ws.Range("G24").Value = 0.002512683913133933 ; $ws.Range("K24").Value = 1.363504477929098 ; $ws.Range("L24").Value = 0.2940851652913352 ; $ws.Range("N24").Value = 2.941854565242934
$ws.Range("B25").Value = 1.288591977884323 ; $ws.Range("D25").Value = 0.1002071453858377 ; $ws.Range("E25").Value = 0.1130814254483257 ; $ws.Range("F25").Value = 2.204868589772261 ; $ws.Range("G25").Value = 0.002521253663993912 ; $ws.Range("K25").Value = 1.182362085673503 ; $ws.Range("L25").Value = 0.2668740032986108 ; $ws.Range("N25").Value = 2.945661970369116
